# Update to user data: add 2016 values for kenya-pop-sex-ratio
# - "Data" sheet: insert a new "2016" row for each district right after its "2014" row
# - "Data-wide-value" sheet: add a "2016" column with the same values
# - "Notes" sheet: update the "Source: 44" note to "Source: 38"

$wb = $excel.ActiveWorkbook

# Map of district_id -> 2016 value
$values2016 = @{
    "d18974" = 93.9
    "d18975" = 103.4
    "d18987" = 91.9
    "d18988" = 82
    "d18976" = 96.2
    "d18955" = 100.7
    "d18965" = 103.8
    "d18968" = 91.4
    "d18956" = 99.8
    "d18991" = 103.1
    "d18989" = 94.4
    "d18977" = 99.7
    "d18943" = 104.6
    "d18949" = 88.4
    "d18944" = 87.8
    "d18969" = 99.1
    "d18970" = 110
    "d18957" = 90.6
    "d18950" = 101.5
    "d18978" = 99.6
    "d18951" = 97.8
    "d18958" = 99.2
    "d18959" = 97.2
    "d18966" = 105.8
    "d18960" = 108.1
    "d18961" = 110.2
    "d18971" = 90
    "d18952" = 103
    "d18946" = 90.4
    "d18964" = 100.5
    "d18979" = 91.8
    "d18980" = 98.1
    "d18981" = 97.2
    "d18962" = 96.3
    "d18972" = 94
    "d18947" = 99.4
    "d18948" = 94.6
    "d18982" = 99.8
    "d18973" = 89.9
    "d18953" = 97.2
    "d18954" = 98.2
    "d18983" = 105.1
    "d18984" = 90.2
    "d18985" = 97.7
    "d18990" = 99.7
    "d18967" = 108.8
    "d18986" = 99
}

# ---------------------------------------------------------------------------
# 1) "Data" sheet - insert a 2016 row after each district's 2014 row
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")
$lastRow = $wsData.UsedRange.Rows.Count

# Walk bottom-up so inserting rows doesn't disturb not-yet-processed rows
for ($r = $lastRow; $r -ge 2; $r--) {
    $year = $wsData.Cells.Item($r, 3).Value2
    if ($year -eq 2014) {
        $districtId = $wsData.Cells.Item($r, 1).Value2
        $districtName = $wsData.Cells.Item($r, 2).Value2
        $newRow = $r + 1
        $wsData.Rows.Item($newRow).Insert()
        $wsData.Cells.Item($newRow, 1).Value2 = $districtId
        $wsData.Cells.Item($newRow, 2).Value2 = $districtName
        $wsData.Cells.Item($newRow, 3).Value2 = 2016
        $wsData.Cells.Item($newRow, 4).Value2 = $values2016[$districtId]
    }
}

# ---------------------------------------------------------------------------
# 2) "Data-wide-value" sheet - add the 2016 column
# ---------------------------------------------------------------------------
$wsWide = $wb.Worksheets.Item("Data-wide-value")

# Force the new header to be stored as text (matching the existing "2013"/"2014"
# text headers) rather than a number, then restore the default cell style so it
# matches its neighbours.
$wsWide.Cells.Item(1, 4).NumberFormat = "@"
$wsWide.Cells.Item(1, 4).Value2 = "2016"
$wsWide.Range("A1").Copy() | Out-Null
$wsWide.Range("D1").PasteSpecial(-4122) | Out-Null

$lastWideRow = $wsWide.UsedRange.Rows.Count
for ($r = 2; $r -le $lastWideRow; $r++) {
    $districtId = $wsWide.Cells.Item($r, 1).Value2
    $wsWide.Cells.Item($r, 4).Value2 = $values2016[$districtId]
}

# ---------------------------------------------------------------------------
# 3) "Notes" sheet - update the source note
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$lastNotesRow = $wsNotes.UsedRange.Rows.Count
for ($r = 1; $r -le $lastNotesRow; $r++) {
    $val = $wsNotes.Cells.Item($r, 1).Value2
    if ($val -eq "Source: 44") {
        $wsNotes.Cells.Item($r, 1).Value2 = "Source: 38"
    }
}
